$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.361.34'
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").Value = '1.658.62'
$ws.Range("E3").Value = '  +0.88%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '235.49'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").Value = '0.4604'
$ws.Range("E7").Value = '  -3.58%  '
$ws.Range("D8").Value = '0.2560'
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("D9").Value = '0.06109'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '1.659.23'
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").Value = '14.52'
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").Value = '4.307'
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("D14").Value = '74.71'
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").Value = '0.5695'
$ws.Range("E15").Value = '  -4.42%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '25.379.06'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '0.000006659'
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").Value = '11.31'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = '1.872.24'
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").Value = '4.390'
$ws.Range("E22").Value = '  +1.43%  '
$ws.Range("D23").Value = '8.589'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '5.183'
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").Value = '133.99'
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").Value = '14.83'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = '1.405'
$ws.Range("E27").Value = '  +1.23%  '
$ws.Range("D28").Value = '1.698'
$ws.Range("E28").Value = '  +3.47%  '
$ws.Range("D29").Value = '103.58'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '3.913'
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("D31").Value = '0.07659'
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").Value = '3.563'
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("D33").Value = '0.04315'
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("D35").Value = '0.5968'
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").Value = '0.9307'
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").Value = '0.9033'
$ws.Range("E37").Value = '  +2.69%  '
$ws.Range("D38").Value = '107.07'
$ws.Range("E38").Value = '  +8.28%  '
$ws.Range("B39").Value = 'PaxDollar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D39").Value = '1.000'
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.311'
$ws.Range("E40").Value = '  -10.37%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '1.813'
$ws.Range("E41").Value = '  +2.62%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01446'
$ws.Range("E42").Value = '  -4.32%  '
$ws.Range("D43").Value = '0.3681'
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("D44").Value = '4.964'
$ws.Range("E44").Value = '  +6.19%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '0.05252'
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1099'
$ws.Range("E46").Value = '  -0.23%  '
$ws.Range("D47").Value = '6.064'
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("D48").Value = '30.13'
$ws.Range("E48").Value = '  +3.76%  '
$ws.Range("D49").Value = '7.540'
$ws.Range("E49").Value = '  +5.84%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '0.9997'
$ws.Range("E51").Value = '  +0.01%  '
